# Update "想去人数" (people interested) counts in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 189
$ws1.Range("F6").Value = 494
$ws1.Range("F9").Value = 40
$ws1.Range("F10").Value = 6583
$ws1.Range("F11").Value = 224
$ws1.Range("F12").Value = 360
$ws1.Range("F13").Value = 2794
$ws1.Range("F14").Value = 172
$ws1.Range("F15").Value = 296
$ws1.Range("F17").Value = 519

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 189
$ws4.Range("F8").Value = 494
$ws4.Range("F11").Value = 40
$ws4.Range("F13").Value = 6583
$ws4.Range("F15").Value = 224
$ws4.Range("F16").Value = 360
$ws4.Range("F17").Value = 2794
$ws4.Range("F18").Value = 172
$ws4.Range("F19").Value = 296
$ws4.Range("F21").Value = 519
